$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=16.58023833333333; H=49.740715; I=0.63541025828417; J=0.63541025828417; K=3; M=4.919927666666666; N=14.759783; O=0.2880104164518184; P=0.2880104164518183; Q=81.57357329609388; R=734.1621596648449; S=0.1830047731061813; T=0.1830047731061812 }
    3  = @{ E=3; G=16.58023833333333; H=49.740715; I=0.63541025828417; J=0.63541025828417; K=3; M=6.683939; N=20.051817; O=0.3912748693382315; P=0.3912748693382315; Q=110.8213016254617; R=997.391714629155; S=0.2486200657863105; T=0.2486200657863105 }
    4  = @{ E=3; G=16.58023833333333; H=49.740715; I=0.63541025828417; J=0.63541025828417; K=3; M=5.478597666666666; N=16.435793; O=0.3207147142099502; P=0.3207147142099501; Q=90.83645504577723; R=817.5280954119951; S=0.2037854193916782; T=0.2037854193916782 }
    5  = @{ E=3; G=8.574149; H=25.722447; I=0.3285901035393414; J=0.3285901035393414; K=3; M=4.919927666666666; N=14.759783; O=0.2880104164518184; P=0.2880104164518183; Q=42.18419288322233; R=379.6577359490009; S=0.09463737256231183; T=0.09463737256231182 }
    6  = @{ E=3; G=8.574149; H=25.722447; I=0.3285901035393414; J=0.3285901035393414; K=3; M=6.683939; N=20.051817; O=0.3912748693382315; P=0.3912748693382315; Q=57.309088892911; R=515.781800036199; S=0.1285690498281918; T=0.1285690498281918 }
    7  = @{ E=3; G=8.574149; H=25.722447; I=0.3285901035393414; J=0.3285901035393414; K=3; M=5.478597666666666; N=16.435793; O=0.3207147142099502; P=0.3207147142099501; Q=46.97431270505233; R=422.768814345471; S=0.1053836811488378; T=0.1053836811488378 }
    8  = @{ E=3; G=0.9393656666666668; H=2.818097; I=0.03599963817648871; J=0.03599963817648871; K=3; M=4.919927666666666; N=14.759783; O=0.2880104164518184; P=0.2880104164518183; Q=4.621611132550111; R=41.594500192951; S=0.01036827078332529; T=0.01036827078332529 }
    9  = @{ E=3; G=0.9393656666666668; H=2.818097; I=0.03599963817648871; J=0.03599963817648871; K=3; M=6.683939; N=20.051817; O=0.3912748693382315; P=0.3912748693382315; Q=6.278662814694334; R=56.50796533224901; S=0.01408575372372923; T=0.01408575372372923 }
    10 = @{ E=3; G=0.9393656666666668; H=2.818097; I=0.03599963817648871; J=0.03599963817648871; K=3; M=5.478597666666666; N=16.435793; O=0.3207147142099502; P=0.3207147142099501; Q=5.146406549546779; R=46.317658945921; S=0.01154561366943419; T=0.01154561366943419 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
